$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the TIPO_VIA value in row 7 (was "CALLE", should be "TRAVESIA")
$ws.Range("C7").Value = "TRAVESIA"

# Update the active selection to C14, matching the saved cursor position
$ws.Range("C14").Select()
